{"js": "const body = context.document.body;\n\n// --- 1 & 2. Title \"# ACTION POINTS Week2 #\" / \"# ACTION POINTS Week3 #\"\n//            -> \"# ACTION POINTS Week4 #\" ---\nconst wk2 = body.search(\"Week2\", { matchCase: true });\nwk2.load(\"items\");\nconst wk3 = body.search(\"Week3\", { matchCase: true });\nwk3.load(\"items\");\nawait context.sync();\n\nwk2.items[0].insertText(\"Week4\", \"Replace\");\nwk3.items[0].insertText(\"Week4\", \"Replace\");\nawait context.sync();\n\n// Paragraph indices are stable up to this point (no paragraphs added/removed yet).\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// --- 3. \"Fix PIP install ( Marco ) shame on you \" -> append \"continue\" ---\nconst p4 = paras.items[3];\np4.insertText(\"continue\", \"End\");\n\n// --- 4. \"Study Operators , Data structure and conditional statements (Iryna)\"\n//        -> \"Complete Exercise 3 file (Iryna)\" ---\nconst p9 = paras.items[8];\np9.insertText(\"Complete Exercise 3 file (Iryna)\", \"Replace\");\n\n// --- 5. \"Complete Exercise 2 file (Iryna)\" -> \"Upload calculator scrip ( Iryna )\" ---\nconst p10 = paras.items[9];\np10.insertText(\"Upload calculator scrip ( Iryna )\", \"Replace\");\n\n// --- 6. \"First Python script Arithmetical Calculator ( Iryna )\" -> \"Study Modules and Loops\" ---\nconst p11 = paras.items[10];\np11.insertText(\"Study Modules and Loops\", \"Replace\");\n\n// --- 7. \"Review Arithmetic Calculator script on git (Marco)\" -> \"Complete debugging Exercise (Iryna)\" ---\nconst p12 = paras.items[11];\np12.insertText(\"Complete debugging Exercise (Iryna)\", \"Replace\");\n\nawait context.sync();\n\n// --- 8. Remove the now-redundant empty paragraph that followed (paraId 0C8E1E40) ---\nconst p13 = paras.items[12];\np13.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Title \"# ACTION POINTS Week2 #\" -> \"# ACTION POINTS Week4 #\" ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"Week2\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"Week4\"\n$find1.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# --- 2. Title \"# ACTION POINTS Week3 #\" -> \"# ACTION POINTS Week4 #\" ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Week3\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Week4\"\n$find2.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# --- 3. \"Fix PIP install ( Marco ) shame on you \" -> append \"continue\" ---\n$p4 = $d.Paragraphs(4)\n$r4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)\n$r4.Collapse(0)\n$r4.InsertAfter(\"continue\")\n\n# --- 4. \"Study Operators , Data structure and conditional statements (Iryna)\"\n#        -> \"Complete Exercise 3 file (Iryna)\" ---\n$p9 = $d.Paragraphs(9)\n$r9 = $d.Range($p9.Range.Start, $p9.Range.End - 1)\n$r9.Text = \"Complete Exercise 3 file (Iryna)\"\n\n# --- 5. \"Complete Exercise 2 file (Iryna)\" -> \"Upload calculator scrip ( Iryna )\" ---\n$p10 = $d.Paragraphs(10)\n$r10 = $d.Range($p10.Range.Start, $p10.Range.End - 1)\n$r10.Text = \"Upload calculator scrip ( Iryna )\"\n\n# --- 6. \"First Python script Arithmetical Calculator ( Iryna )\" -> \"Study Modules and Loops\" ---\n$p11 = $d.Paragraphs(11)\n$r11 = $d.Range($p11.Range.Start, $p11.Range.End - 1)\n$r11.Text = \"Study Modules and Loops\"\n\n# --- 7. \"Review Arithmetic Calculator script on git (Marco)\" -> \"Complete debugging Exercise (Iryna)\" ---\n$p12 = $d.Paragraphs(12)\n$r12 = $d.Range($p12.Range.Start, $p12.Range.End - 1)\n$r12.Text = \"Complete debugging Exercise (Iryna)\"\n\n# --- 8. Remove the now-redundant empty paragraph that followed (paraId 0C8E1E40) ---\n$p13 = $d.Paragraphs(13)\n$p13.Range.Delete()\n"}
